# Apply cryptos list update (price/volume refresh + two row swaps)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.188.09"
$ws.Range("E2").Value = "  +0.35%  "
$ws.Range("D3").Value = "1.600.12"
$ws.Range("E3").Value = "  -0.11%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D6").Value = "'303.18"
$ws.Range("E6").Value = "  +0.40%  "
$ws.Range("D7").Value = "'0.3780"
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "'52.18"
$ws.Range("E8").Value = "  +4.76%  "
$ws.Range("E9").Value = "  -1.18%  "
$ws.Range("E10").Value = "  -0.58%  "
$ws.Range("D11").Value = "'1.001"
$ws.Range("E11").Value = "  +0.01%  "
$ws.Range("D12").Value = "'0.08116"
$ws.Range("E12").Value = "  -0.58%  "
$ws.Range("D13").Value = "'22.66"
$ws.Range("E13").Value = "  -2.17%  "
$ws.Range("E14").Value = "  -0.43%  "
$ws.Range("D15").Value = "'7.399"
$ws.Range("E15").Value = "  +0.06%  "
$ws.Range("D16").Value = "'0.00001245"
$ws.Range("E16").Value = "  -1.29%  "
$ws.Range("D17").Value = "1.600.85"
$ws.Range("E17").Value = "  -0.26%  "
$ws.Range("E18").Value = "  +2.60%  "
$ws.Range("D19").Value = "'0.06890"
$ws.Range("E19").Value = "  +0.61%  "
$ws.Range("D20").Value = "'18.06"
$ws.Range("E20").Value = "  -1.73%  "
$ws.Range("D21").Value = "'6.542"
$ws.Range("E21").Value = "  -0.57%  "
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("E23").Value = "  -0.61%  "
$ws.Range("D24").Value = "23.183.20"
$ws.Range("E24").Value = "  +0.35%  "
$ws.Range("B25").Value = "LidoDAOToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D25").Value = "'3.002"
$ws.Range("E25").Value = "  +10.34%  "
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").Value = "'2.383"
$ws.Range("E26").Value = "  +1.86%  "
$ws.Range("E27").Value = "  +0.35%  "
$ws.Range("D28").Value = "'149.52"
$ws.Range("E28").Value = "  -0.73%  "
$ws.Range("D29").Value = "'5.251"
$ws.Range("E29").Value = "  -0.48%  "
$ws.Range("D30").Value = "'133.81"
$ws.Range("D31").Value = "'2.384"
$ws.Range("E31").Value = "  -0.74%  "
$ws.Range("D32").Value = "'6.796"
$ws.Range("E32").Value = "  -0.79%  "
$ws.Range("D33").Value = "1.779.78"
$ws.Range("E33").Value = "  -0.54%  "
$ws.Range("D34").Value = "'0.9657"
$ws.Range("E34").Value = "  +0.53%  "
$ws.Range("D35").Value = "'0.07489"
$ws.Range("E35").Value = "  -2.59%  "
$ws.Range("D36").Value = "'10.27"
$ws.Range("E36").Value = "  +1.67%  "
$ws.Range("D37").Value = "'0.02705"
$ws.Range("E37").Value = "  -1.07%  "
$ws.Range("D38").Value = "'0.2504"
$ws.Range("E38").Value = "  -2.02%  "
$ws.Range("D39").Value = "'0.08802"
$ws.Range("E39").Value = "  -1.28%  "
$ws.Range("D40").Value = "'6.080"
$ws.Range("E40").Value = "  -3.61%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "'0.7098"
$ws.Range("E41").Value = "  -0.12%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "'1.361"
$ws.Range("E42").Value = "  -0.67%  "
$ws.Range("E43").Value = "  -1.76%  "
$ws.Range("D44").Value = "'15.46"
$ws.Range("E44").Value = "  +0.61%  "
$ws.Range("D45").Value = "'0.6520"
$ws.Range("E45").Value = "  -1.71%  "
$ws.Range("D47").Value = "'4.012"
$ws.Range("E47").Value = "  +0.53%  "
$ws.Range("D48").Value = "'132.03"
$ws.Range("E48").Value = "  -0.27%  "
$ws.Range("D49").Value = "'0.07961"
$ws.Range("E49").Value = "  +0.23%  "
$ws.Range("D50").Value = "'1.200"
$ws.Range("E50").Value = "  -0.57%  "
$ws.Range("D51").Value = "'1.214"
$ws.Range("E51").Value = "  +1.78%  "
